$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the style used by the other headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the new Save column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
